# Training file update: new response-distribution data entered for
# questions in rows 16-19, and the "number correct weighting" (D24)
# switched on (0 -> 1). Also resets the view (selection/scroll position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Row 16 ---
$ws.Range("D16").Value = 0.692
$ws.Range("E16:P16").Value = 0.02567

# --- Row 17 ---
$ws.Range("D17").Value = 0.764
$ws.Range("E17:P17").Value = 0.01967

# --- Row 18 ---
$ws.Range("D18").Value = 0.824
$ws.Range("E18:P18").Value = 0.01467

# --- Row 19 ---
$ws.Range("D19").Value = 0.848
$ws.Range("E19:P19").Value = 0.01267

# --- Row 24 (enable the IRT/weighting factor) ---
$ws.Range("D24").Value = 1

# --- Reset the view: clear the scrolled/selected state saved in the file ---
$ws.Activate()
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
